$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[name=`"'Deathless Black Snake'`"]  Well done.`n"
$ws.Range("C3").Value = "[name=`"'Deathless Black Snake'`"]  You’ve done well...`n"
$ws.Range("C4").Value = "[name=`"'Deathless Black Snake'`"]  What is this? The dragonslaying sword couldn’t cut me, the sword of the Lord of Fiends’ couldn’t pierce me, so why is... my Arts... My control is fading away?`n"
$ws.Range("C31").Value = "[name=`"Talulah?`"] Have you ever seen someone cut his throat in front of you, 'sister'?`n"
$ws.Range("C41").Value = "[name=`"Talulah?`"] I am merely 'educating' her.`n"
$ws.Range("C42").Value = "[name=`"Talulah?`"]  My failure stems from my enduring persistence. When you encounter a person without persistence, you will come to realize how much is contained within the word 'suffering'.`n"
$ws.Range("C45").Value = "[name=`"Talulah?`"] Even I myself am just part of the road she steps upon. If I had succeeded in killing you, 'sister'... It’d be a lot simpler.`n"
$ws.Range("C49").Value = "[name=`"Ch’en`"]  You have the nerve to say that?! After sacrificing your own 'daughter'? He’s nothing like you! `n"
$ws.Range("C70").Value = "This is definitely not the first time the white-haired Draco attempted to resist her 'father'.`n"
$ws.Range("C81").Value = "[Subtitle(text=`"The different between me and you is that your 'love'... is just sacrifice.`", x=200, y=360, alignment=`"left`", size=24, delay=0.04, width=1280)]`n"
$ws.Range("C92").Value = "[name=`"Talulah`"]  'I taught you all of this?' No, Kashchey. This world, that snowfield, those people chasing after the sunlight... You will never understand the things they taught me.`n"
$ws.Range("D101").Value = "[name=`"아미야`"]  체, 첸 씨? 그런 말은…… 피디아 족 오퍼레이터들한텐 하시면 안 돼요!`n"
$ws.Range("C117").Value = "[name=`"Talulah?`"]  'I will be there even at the edge of the world.' `n"
$ws.Range("D129").Value = "[name=`"W`"]  아, 참! 그리고…… 날 안 믿는다고 해도 상관없어. 나도 너 안 믿거든. 하지만 테레시아의 후계자가 이런 순간에 멍청한 짓을 하진 않을 거라고 생각하니까……`n"
$ws.Range("D166").Value = "[name=`"W`"]  이번 한 번만 부탁할게!`n"
$ws.Range("C177").Value = "[name=`"W`"]  What the hell does it mean 'insufficient clearance'?`n"
$ws.Range("C289").Value = "[name=`"Talulah`"]  'Good food, great wine and amazing scenery; Pretty women, good manners and great company'. Isn’t it like those weird novels you used to read?`n"
